# "unify the conception of DataNode, DataTable, Entity."
#
# The two worksheets are renamed to the new, unified naming scheme, and the
# workbook is left the way it was last saved by the author: with the
# DataTable sheet active and a cell further down it selected. A couple of
# header-row heights on the DataNode sheet are also adjusted.

$wb = $excel.ActiveWorkbook

$wsNode  = $wb.Worksheets.Item(1)   # was "Property1"
$wsTable = $wb.Worksheets.Item(2)   # was "Record_Station"

$wsNode.Name  = "DataNode"
$wsTable.Name = "DataTable"

# Header row wraps a little taller; the column-description row a little
# shorter.
$wsNode.Rows.Item(1).RowHeight = 27
$wsNode.Rows.Item(8).RowHeight = 54

# Workbook was last saved on the DataTable sheet, with H32 selected.
$wsTable.Activate()
$null = $wsTable.Range("H32").Select()
